$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the worksheet "Clara" -> "Maia & Clara"
#    (the _xlnm._FilterDatabase defined name is re-pointed automatically by
#    the host when the sheet is renamed)
# ---------------------------------------------------------------------------
$ws.Name = "Maia & Clara"

# ---------------------------------------------------------------------------
# 2. Mark a handful of already-present wishes as reserved/bought (column E)
#    and hide those rows, same as ticking them off in the filtered sheet.
# ---------------------------------------------------------------------------
$reservedRows = @(4, 5, 9, 10, 12, 13, 17)
foreach ($r in $reservedRows) {
    $ws.Range("E$r").Value = "Y"
    $ws.Rows.Item($r).Hidden = $true
}

# ---------------------------------------------------------------------------
# 3. Append the newly-wished-for items as rows 24-31
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = "Yoshi''s Crafted World"
$ws.Range("B24").Value = "https://static01.galaxus.com/productimages/8/3/7/3/1/4/5/9/8/4/4/5/6/6/0/2/4/4/5/2ba1a351-3eaf-4c4a-9efc-f1e6a6c8f996_cropped.jpg_480.avif"
$ws.Range("C24").Value = "https://www.digitec.ch/en/s1/product/nintendo-yoshis-crafted-world-switch-switch-lite-switch-oled-en-video-games-20454417"
$ws.Range("D24").Value = "49 CHF"

$ws.Range("A25").Value = "Captain Toad: Treasure Tracker"
$ws.Range("B25").Value = "https://www.wog.ch/nas/cover_large/sw/sw_captaintoadtreasuretracker.jpg"
$ws.Range("C25").Value = "https://www.wog.ch/de/index.cfm/details/product/63694-Captain-Toad-Treasure-Tracker"
$ws.Range("D25").Value = "37.90 CHF"

$ws.Range("A26").Value = "Sticla de apa div. modele (2)"
$ws.Range("B26").Value = "https://static01.galaxus.com/productimages/5/2/5/8/5/2/2/9/9/7/8/8/3/6/1/7/3/4/2/49014e07-1c9f-43f9-ad17-6ec479dba044_cropped.jpg_2880.avif"
$ws.Range("C26").Value = "https://www.galaxus.ch/en/s3/product/24-bottles-urban-025-l-water-bottles-thermos-flasks-21389032"
$ws.Range("D26").Value = "16.90 CHF"

$ws.Range("A27").Value = "Sticla de apa div. modele (2)"
$ws.Range("B27").Value = "https://static01.galaxus.com/productimages/3/7/2/0/5/9/2/6/6/3/0/4/1/2/7/1/9/3/8/d1fc0547-b35f-4453-ae72-b6dae88e9d7d.png_480.avif"
$ws.Range("C27").Value = "https://www.galaxus.ch/en/s3/product/sigg-miracle-wmb-040-l-water-bottles-thermos-flasks-53291495"
$ws.Range("D27").Value = "23.90 CHF"

$ws.Range("A28").Value = "Furzipups der Knatterdrache"
$ws.Range("B28").Value = "https://image.smythstoys.com/zoom/246869.webp"
$ws.Range("C28").Value = "https://www.smythstoys.com/ch/de-ch/spielzeug/brettspiele-und-gesellschaftsspiele/kinderspiele/furzipups-der-knatterdrache/p/246869"
$ws.Range("D28").Value = "9.95 CHF"

$ws.Range("A29").Value = "Schnappt Hubi!"
$ws.Range("B29").Value = "https://image.smythstoys.com/zoom/8000281_5.webp"
$ws.Range("C29").Value = "https://www.smythstoys.com/ch/de-ch/spielzeug/brettspiele-und-gesellschaftsspiele/kinderspiele/schnappt-hubi-spiel-mit-dem-frechen-gespenst/p/8000281"
$ws.Range("D29").Value = "34.95 CHF"

$ws.Range("A30").Value = "Paletti Spaghetti"
$ws.Range("B30").Value = "https://image.smythstoys.com/zoom/209042.webp"
$ws.Range("C30").Value = "https://www.smythstoys.com/ch/de-ch/spielzeug/brettspiele-und-gesellschaftsspiele/kinderspiele/paletti-spaghetti/p/209042"
$ws.Range("D30").Value = "19.95 CHF"

$ws.Range("A31").Value = "Sticla de apa div. modele (2)"
$ws.Range("B31").Value = "https://static01.galaxus.com/productimages/1/0/1/8/5/4/1/0/4/7/6/6/3/7/5/6/5/9/5/a9356c14-5bd5-48c2-a32c-c461302d2146.jpg_480.avif"
$ws.Range("C31").Value = "https://www.galaxus.ch/en/s3/product/camelbak-thrive-flip-straw-kids-vi-040-l-water-bottles-thermos-flasks-52764026"
$ws.Range("D31").Value = "27.90 CHF"

# ---------------------------------------------------------------------------
# 4. Re-apply the AutoFilter so it spans the grown table (A1:E31) while
#    keeping the same "blank" filter on column E (5th column).
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:E31").AutoFilter(5, @(""), 7)

# ---------------------------------------------------------------------------
# 5. The hidden _xlnm._FilterDatabase name isn't re-synced by the AutoFilter
#    resize above, so point it at the new range/sheet name explicitly.
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Maia & Clara'!`$A`$1:`$E`$31"
    }
}

# ---------------------------------------------------------------------------
# 6. Put the selection where the author last left it (B28) for parity with
#    the saved workbook view.
# ---------------------------------------------------------------------------
$ws.Range("B28").Select()
